# Saldo.xlsx edit: remove a handful of account rows, update a couple of
# balances, and relocate the LUCYENE row further down the list (with its
# own balance update), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase A: delete whole rows (bottom-most row index first so the
#     remaining row numbers we still need don't shift under us) ---
$ws.Rows.Item(15).Delete()   # 004550415 DIOGO       2000
$ws.Rows.Item(14).Delete()   # 004805333 ROSANA       3802.6
$ws.Rows.Item(13).Delete()   # 004265173 JULIA        5306.54
$ws.Rows.Item(5).Delete()    # 004480970 ALBERTO      67645.49
$ws.Rows.Item(2).Delete()    # 004212438 KENIA        290404.58

# --- Phase B: straightforward balance updates ---
$ws.Cells.Item(3, 3).Value = 50030.14    # 004459461 INTERLAGOS
$ws.Cells.Item(11, 3).Value = 1405.27    # 004467884 ANA

# --- Phase C: relocate the LUCYENE row (currently row 6) to just before
#     GUSTAVO (currently row 21), preserving its original cell formatting
#     by cutting/pasting the row rather than retyping the account number ---
$ws.Rows.Item(21).Insert()
$ws.Range("A6:C6").Cut($ws.Range("A21:C21"))
$ws.Rows.Item(6).Delete()

# --- Phase D: LUCYENE's own balance update, now at row 20 ---
$ws.Cells.Item(20, 3).Value = 209.92
